$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = 'D2'; Value = '28.365.48' }
    @{ Cell = 'E2'; Value = '  +0.17%  ' }
    @{ Cell = 'D3'; Value = '1.866.03' }
    @{ Cell = 'E3'; Value = '  -0.22%  ' }
    @{ Cell = 'E4'; Value = '  +0.13%  ' }
    @{ Cell = 'D5'; Value = '330.78' }
    @{ Cell = 'E5'; Value = '  -2.66%  ' }
    @{ Cell = 'D6'; Value = '1.002' }
    @{ Cell = 'E6'; Value = '  +0.09%  ' }
    @{ Cell = 'D7'; Value = '0.4623' }
    @{ Cell = 'E7'; Value = '  -1.77%  ' }
    @{ Cell = 'D8'; Value = '0.4004' }
    @{ Cell = 'E8'; Value = '  +1.58%  ' }
    @{ Cell = 'D9'; Value = '47.79' }
    @{ Cell = 'E9'; Value = '  +1.07%  ' }
    @{ Cell = 'D10'; Value = '0.07859' }
    @{ Cell = 'E10'; Value = '  -1.71%  ' }
    @{ Cell = 'E11'; Value = '  -2.33%  ' }
    @{ Cell = 'D12'; Value = '21.27' }
    @{ Cell = 'E12'; Value = '  -2.99%  ' }
    @{ Cell = 'D13'; Value = '1.864.14' }
    @{ Cell = 'E13'; Value = '  -1.06%  ' }
    @{ Cell = 'D14'; Value = '5.848' }
    @{ Cell = 'E14'; Value = '  -2.63%  ' }
    @{ Cell = 'D15'; Value = '6.993' }
    @{ Cell = 'E15'; Value = '  -4.01%  ' }
    @{ Cell = 'E16'; Value = '  +0.04%  ' }
    @{ Cell = 'D17'; Value = '88.17' }
    @{ Cell = 'E17'; Value = '  -3.33%  ' }
    @{ Cell = 'D18'; Value = '0.06548' }
    @{ Cell = 'E18'; Value = '  -0.71%  ' }
    @{ Cell = 'D19'; Value = '0.00001018' }
    @{ Cell = 'E19'; Value = '  -2.35%  ' }
    @{ Cell = 'D20'; Value = '17.20' }
    @{ Cell = 'E20'; Value = '  -2.77%  ' }
    @{ Cell = 'D21'; Value = '0.9998' }
    @{ Cell = 'E21'; Value = '  -0.13%  ' }
    @{ Cell = 'D22'; Value = '28.348.49' }
    @{ Cell = 'E22'; Value = '  +0.08%  ' }
    @{ Cell = 'E23'; Value = '  -2.09%  ' }
    @{ Cell = 'D24'; Value = '10.87' }
    @{ Cell = 'E24'; Value = '  -1.76%  ' }
    @{ Cell = 'D25'; Value = '2.236' }
    @{ Cell = 'E25'; Value = '  -2.42%  ' }
    @{ Cell = 'D26'; Value = '2.089.68' }
    @{ Cell = 'E26'; Value = '  -0.73%  ' }
    @{ Cell = 'D27'; Value = '157.34' }
    @{ Cell = 'E27'; Value = '  -1.57%  ' }
    @{ Cell = 'D28'; Value = '19.35' }
    @{ Cell = 'E28'; Value = '  -2.45%  ' }
    @{ Cell = 'E29'; Value = '  -4.27%  ' }
    @{ Cell = 'D30'; Value = '5.293' }
    @{ Cell = 'E30'; Value = '  -3.75%  ' }
    @{ Cell = 'D31'; Value = '117.63' }
    @{ Cell = 'E31'; Value = '  -2.31%  ' }
    @{ Cell = 'D32'; Value = '0.9592' }
    @{ Cell = 'E32'; Value = '  -1.79%  ' }
    @{ Cell = 'D33'; Value = '0.09346' }
    @{ Cell = 'E33'; Value = '  -1.60%  ' }
    @{ Cell = 'D34'; Value = '3.583' }
    @{ Cell = 'E34'; Value = '  -0.26%  ' }
    @{ Cell = 'D35'; Value = '1.384' }
    @{ Cell = 'E35'; Value = '  +0.19%  ' }
    @{ Cell = 'D36'; Value = '5.249' }
    @{ Cell = 'E36'; Value = '  -2.03%  ' }
    @{ Cell = 'D37'; Value = '0.06036' }
    @{ Cell = 'E37'; Value = '  -0.96%  ' }
    @{ Cell = 'D38'; Value = '0.02202' }
    @{ Cell = 'E38'; Value = '  -3.40%  ' }
    @{ Cell = 'D39'; Value = '8.265' }
    @{ Cell = 'E39'; Value = '  -2.29%  ' }
    @{ Cell = 'E40'; Value = '  -1.69%  ' }
    @{ Cell = 'D41'; Value = '1.001' }
    @{ Cell = 'E41'; Value = '  +0.03%  ' }
    @{ Cell = 'D42'; Value = '0.5757' }
    @{ Cell = 'B43'; Value = 'Algorand' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Cell = 'D43'; Value = '0.1810' }
    @{ Cell = 'E43'; Value = '  -3.85%  ' }
    @{ Cell = 'B44'; Value = 'Aptos' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'D44'; Value = '10.06' }
    @{ Cell = 'E44'; Value = '  -3.27%  ' }
    @{ Cell = 'D45'; Value = '1.269' }
    @{ Cell = 'E45'; Value = '  -2.87%  ' }
    @{ Cell = 'D46'; Value = '2.289' }
    @{ Cell = 'E46'; Value = '  +13.08%  ' }
    @{ Cell = 'B47'; Value = 'Decentraland' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D47'; Value = '0.5431' }
    @{ Cell = 'E47'; Value = '  -3.54%  ' }
    @{ Cell = 'B48'; Value = 'EnergySwap' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D48'; Value = '11.84' }
    @{ Cell = 'E48'; Value = '  -2.18%  ' }
    @{ Cell = 'D49'; Value = '0.07118' }
    @{ Cell = 'E49'; Value = '  +2.89%  ' }
    @{ Cell = 'D50'; Value = '1.889' }
    @{ Cell = 'E50'; Value = '  -4.02%  ' }
    @{ Cell = 'D51'; Value = '111.28' }
    @{ Cell = 'E51'; Value = '  -0.07%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"
